$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 4 (event 14552522) with the match result: "Fallo" (-1 profit)
$ws.Range("G4").Value = "Fallo"
$ws.Range("H4").Value = -1
